$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 15
$ws.Range("H15").Value = 770937.25
$ws.Range("I15").Value = 770937.25
$ws.Range("K15").Value = 2312811.75
$ws.Range("M15").Value = -2312642.75

# row 32
$ws.Range("H32").Value = 590
$ws.Range("J32").Value = 590
$ws.Range("L32").Value = 590
$ws.Range("N32").Value = -1242

# row 51
$ws.Range("H51").Value = 3499.1667
$ws.Range("J51").Value = 3599
$ws.Range("L51").Value = 3599
$ws.Range("N51").Value = -4567

# row 70
$ws.Range("H70").Value = 3035.4333
$ws.Range("J70").Value = 5928.5
$ws.Range("L70").Value = 17785.5
$ws.Range("N70").Value = -18325.5

# row 73
$ws.Range("H73").Value = 3035.4333
$ws.Range("J73").Value = 5928.5
$ws.Range("L73").Value = 17785.5
$ws.Range("N73").Value = -19657.5

# row 88
$ws.Range("H88").Value = 12087.444
$ws.Range("J88").Value = 12473.375
$ws.Range("L88").Value = 12473.375
$ws.Range("N88").Value = -13285.375

# row 91
$ws.Range("H91").Value = 12087.444
$ws.Range("J91").Value = 12473.375
$ws.Range("L91").Value = 12473.375
$ws.Range("N91").Value = -15281.375

# row 96
$ws.Range("H96").Value = 1837.75
$ws.Range("I96").Value = 1837.75
$ws.Range("K96").Value = 5513.25
$ws.Range("M96").Value = -4140.25

# row 106
$ws.Range("H106").Value = 1730.579
$ws.Range("I106").Value = 1748.9445
$ws.Range("K106").Value = 1748.9445
$ws.Range("M106").Value = -1117.9445

# row 113
$ws.Range("H113").Value = 4391
$ws.Range("I113").Value = 3966.6667
$ws.Range("K113").Value = 3966.6667
$ws.Range("M113").Value = -712.6667000000002

$ws = $wb.Worksheets.Item("ARM")
# row 61
$ws.Range("H61").Value = 41668044
$ws.Range("I61").Value = 50001130
$ws.Range("K61").Value = 50001130
$ws.Range("M61").Value = -50000918

# row 74
$ws.Range("H74").Value = 26320368
$ws.Range("I74").Value = 30307530
$ws.Range("J74").Value = 5103.8
$ws.Range("K74").Value = 30307530
$ws.Range("L74").Value = 5103.8
$ws.Range("M74").Value = -30306656
$ws.Range("N74").Value = -6851.8

# row 77
$ws.Range("H77").Value = 26320368
$ws.Range("I77").Value = 30307530
$ws.Range("J77").Value = 5103.8
$ws.Range("K77").Value = 151537650
$ws.Range("L77").Value = 25519
$ws.Range("M77").Value = -151533282
$ws.Range("N77").Value = -34255

# row 132
$ws.Range("H132").Value = 3708341.8
$ws.Range("I132").Value = 4170921.8
$ws.Range("K132").Value = 12512765.4
$ws.Range("M132").Value = -12510235.4

# row 136
$ws.Range("H136").Value = 41668044
$ws.Range("I136").Value = 50001130
$ws.Range("K136").Value = 150003390
$ws.Range("M136").Value = -150000840

# row 140
$ws.Range("H140").Value = 206284.67
$ws.Range("J140").Value = 206284.67
$ws.Range("L140").Value = 206284.67
$ws.Range("N140").Value = -216644.67

$ws = $wb.Worksheets.Item("BSM")
# row 22
$ws.Range("H22").Value = 3249.8
$ws.Range("I22").Value = 4700.6
$ws.Range("K22").Value = 4700.6
$ws.Range("M22").Value = -4527.6

# row 140
$ws.Range("H140").Value = 98725
$ws.Range("J140").Value = 98725
$ws.Range("L140").Value = 98725
$ws.Range("N140").Value = -109085

$ws = $wb.Worksheets.Item("CRP")
# row 22
$ws.Range("H22").Value = 745
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()

# row 32
$ws.Range("H32").Value = 6957.25
$ws.Range("I32").Value = 6957.25
$ws.Range("K32").Value = 6957.25
$ws.Range("M32").Value = -6641.25

# row 86
$ws.Range("H86").Value = 11427.261
$ws.Range("I86").Value = 8969.25
$ws.Range("J86").Value = 14108.728
$ws.Range("K86").Value = 8969.25
$ws.Range("L86").Value = 14108.728
$ws.Range("M86").Value = -7846.25
$ws.Range("N86").Value = -16354.728

# row 89
$ws.Range("H89").Value = 11427.261
$ws.Range("I89").Value = 8969.25
$ws.Range("J89").Value = 14108.728
$ws.Range("K89").Value = 44846.25
$ws.Range("L89").Value = 70543.64
$ws.Range("M89").Value = -39230.25
$ws.Range("N89").Value = -81775.64

# row 99
$ws.Range("H99").Value = 16174.75
$ws.Range("I99").Value = 16174.75
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 16174.75
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -14676.75
$ws.Range("N99").ClearContents()

# row 126
$ws.Range("H126").Value = 16174.75
$ws.Range("I126").Value = 16174.75
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 48524.25
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -46054.25
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# row 80
$ws.Range("H80").Value = 3584.6
$ws.Range("I80").Value = 3584.6
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 10753.8
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -9817.799999999999
$ws.Range("N80").ClearContents()

# row 83
$ws.Range("H83").Value = 3584.6
$ws.Range("I83").Value = 3584.6
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 32261.4
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -27581.4
$ws.Range("N83").ClearContents()

# row 86
$ws.Range("H86").Value = 725.94116
$ws.Range("I86").Value = 640.2857
$ws.Range("J86").Value = 785.9
$ws.Range("K86").Value = 1920.8571
$ws.Range("L86").Value = 2357.7
$ws.Range("M86").Value = -734.8571000000002
$ws.Range("N86").Value = -4729.7

# row 89
$ws.Range("H89").Value = 725.94116
$ws.Range("I89").Value = 640.2857
$ws.Range("J89").Value = 785.9
$ws.Range("K89").Value = 5762.571300000001
$ws.Range("L89").Value = 7073.099999999999
$ws.Range("M89").Value = 165.4286999999995
$ws.Range("N89").Value = -18929.1

$ws = $wb.Worksheets.Item("GSM")
# row 15
$ws.Range("H15").Value = 44138
$ws.Range("J15").Value = 44138
$ws.Range("L15").Value = 44138
$ws.Range("N15").Value = -44714

# row 81
$ws.Range("H81").Value = 44138
$ws.Range("J81").Value = 44138
$ws.Range("L81").Value = 44138
$ws.Range("N81").Value = -46134

# row 84
$ws.Range("H84").Value = 44138
$ws.Range("J84").Value = 44138
$ws.Range("L84").Value = 132414
$ws.Range("N84").Value = -142398

# row 132
$ws.Range("H132").Value = 5439806.5
$ws.Range("I132").Value = 6583534
$ws.Range("K132").Value = 19750602
$ws.Range("M132").Value = -19748072

# row 134
$ws.Range("H134").Value = 50000
$ws.Range("J134").Value = 50000
$ws.Range("L134").Value = 150000
$ws.Range("N134").Value = -155070

# row 136
$ws.Range("H136").Value = 45000
$ws.Range("J136").Value = 45000
$ws.Range("L136").Value = 135000
$ws.Range("N136").Value = -140100

$ws = $wb.Worksheets.Item("LTW")
# row 22
$ws.Range("H22").Value = 2844.4443
$ws.Range("I22").Value = 3164.2856
$ws.Range("J22").Value = 1725
$ws.Range("K22").Value = 3164.2856
$ws.Range("L22").Value = 1725
$ws.Range("M22").Value = -2869.2856
$ws.Range("N22").Value = -2315

# row 27
$ws.Range("H27").Value = 2844.4443
$ws.Range("I27").Value = 3164.2856
$ws.Range("J27").Value = 1725
$ws.Range("K27").Value = 3164.2856
$ws.Range("L27").Value = 1725
$ws.Range("M27").Value = -3057.2856
$ws.Range("N27").Value = -1939

# row 118
$ws.Range("H118").Value = 39499.5
$ws.Range("J118").Value = 39499.5
$ws.Range("L118").Value = 39499.5
$ws.Range("N118").Value = -42813.5

$ws = $wb.Worksheets.Item("WVR")
# row 62
$ws.Range("H62").Value = 4250
$ws.Range("I62").Value = 4250
$ws.Range("K62").Value = 4250
$ws.Range("M62").Value = -3626

# row 65
$ws.Range("H65").Value = 4250
$ws.Range("I65").Value = 4250
$ws.Range("K65").Value = 21250
$ws.Range("M65").Value = -18130

# row 100
$ws.Range("H100").Value = 879.3461
$ws.Range("I100").Value = 820.375
$ws.Range("J100").Value = 1587
$ws.Range("K100").Value = 1640.75
$ws.Range("L100").Value = 3174
$ws.Range("M100").Value = -1099.75
$ws.Range("N100").Value = -4256

# row 107
$ws.Range("H107").Value = 609.8182
$ws.Range("I107").Value = 585.8
$ws.Range("J107").Value = 850
$ws.Range("K107").Value = 1757.4
$ws.Range("L107").Value = 2550
$ws.Range("M107").Value = 162.6000000000001
$ws.Range("N107").Value = -6390

# row 122
$ws.Range("H122").Value = 3249.75
$ws.Range("I122").Value = 2500
$ws.Range("K122").Value = 7500
$ws.Range("M122").Value = -5050
